$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 4 data
$ws.Cells.Item(4, 1).Value = 9102000
$ws.Cells.Item(4, 2).Value = 1
$ws.Cells.Item(4, 3).Value = 200
$ws.Cells.Item(4, 4).Value = 0
$ws.Cells.Item(4, 5).Value = 0
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 0
$ws.Cells.Item(4, 8).Value = 1
$ws.Cells.Item(4, 9).Value = "サファイアSP"

# Adjust row heights for rows 2 and 3
$ws.Rows.Item(2).RowHeight = 16
$ws.Rows.Item(3).RowHeight = 16

# Select A4 to match final selection state
$ws.Range("A4").Select()
